# Schedule Metrics Tracking — Iteration 5 update
# Row 8 corresponds to "Iteration 5" in Table1 (B3:K18 on "Schedule Metrics Tracking").
# The team's actual start slipped by 28 days versus the original plan, and the actual
# duration came in at 15 days instead of the estimated 1 (placeholder) day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule Metrics Tracking")

# Planned Start Date (C8) is corrected from the formula-derived date to a literal date,
# which also cascades (via the existing Planned End/Start Date formulas) through every
# later iteration row (9-18).
$ws.Range("C8").Value = 41988

# Actual Start Date / Actual End Date (G8/H8) get filled in now that Iteration 5 finished.
$ws.Range("G8").Value = 41988
$ws.Range("H8").Formula = "=G8+13+1"

# Actual Duration (Days) updates from the placeholder 1 to the real 15 days taken.
$ws.Range("I8").Value = 15

# Bring the external "Schedule Metrics" workbook link (used by the Gantt charts'
# filtered category titles) back into the link table, matching the source workbook.
$ws.Range("M200").Formula = "='[Schedule Metrics.xlsx]Schedule Metrics'!B7"
$ws.Range("M200").ClearContents()

# Restore the view state left behind by the editor (scrolled down, zoomed to 70%,
# with I8 selected after entering the actual duration).
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("I8").Select()
$excel.ActiveWindow.ScrollRow = 31
